$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles 24/25, borders, wrap text) from the last existing row
$sourceRow = 1028
$ws.Range("A" + $sourceRow + ":G" + $sourceRow).Copy() | Out-Null
$ws.Range("A1029:G1040").PasteSpecial(-4122) | Out-Null

# Row 1029
$ws.Range("A1029").Value = "TN"
$ws.Range("B1029").Value = "ELUR03_CBE_P40"
$ws.Range("C1029").Value = "28-Dec-2025 12:55 PM"
$ws.Range("D1029").Value = "FAIL"
$ws.Range("E1029").Value = "1. SCG addition after VoLTE call released`n2. Video Streaming  (ms)"
$ws.Range("F1029").Value = "1. Static VoLTE MO`n2. Static Yotube Streaming"
$ws.Range("G1029").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n2. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows(1029).RowHeight = 145

# Row 1030
$ws.Range("A1030").Value = "CH"
$ws.Range("B1030").Value = "VIJNA8_CHN_P41"
$ws.Range("C1030").Value = "28-Dec-2025 12:57 PM"
$ws.Range("D1030").Value = "FAIL"
$ws.Range("E1030").Value = "1. RACH setup`n2. MO Call (pass/fail)`n3. MT Call (pass/fail)`n4. Downlink Peak MCS - 5G`n5. Ping/Round trip time(ms)`n6. UE Steering (Connected) : Non anchor/anchor to preferred anchor`n7. QCI Verification`n8. Web Browsing - Top 10 Websites - Web page load time (ms)"
$ws.Range("F1030").Value = "1. Static All`n2. Static VoLTE MO`n3. Static VoLTE MT`n4. Static DL`n5. Static Ping`n6. Static DL`n7. Static All`n8. Static Browsing (10 sites)"
$ws.Range("G1030").Value = "1. If DT Tool is  TEMS Pocket, verify the Static ATDT . The NR RACH Attempts should be equal to NR RACH Success; kindly exclude the logs where NR RACH has failed.`nIf DT Tool is  AZQ, Validate Static All and ensure NR RACH Attempts match NR RACH Success. Please exclude the logs with NR RACH failures and redo the test accordingly.`n2. VoLTE Long Call MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO calls.`n3. VoLTE Long Call MT – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MT calls.`n4. Peak MCS is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n5. Ping is not meeting the acceptance criteria. The average ping value across all logfiles should be less than 50 ms. Kindly exclude the logfile where the average value exceeds 50 ms and redo the test.”`n6. Kindly verify that the UE in connected mode is switching from the preferred anchor layer to NR as expected.`n7. If the DT tool is TEMS, verify the QCI combined value in the DL drive. If the DT tool is AZQ, verify it in Static All. The QCI combined value must meet the acceptance criteria.`n8. While running the scripts, kindly verify in the Events tab that at least 10 websites are browsing. If not, exclude the logfile and create a new one."
$ws.Rows(1030).RowHeight = 377

# Row 1031
$ws.Range("A1031").Value = "TN"
$ws.Range("B1031").Value = "CB1217_CBE_P40"
$ws.Range("C1031").Value = "28-Dec-2025 5:56 PM"
$ws.Range("D1031").Value = "FAIL"
$ws.Range("E1031").Value = "1. Video Streaming  (ms)"
$ws.Range("F1031").Value = "1. Static Yotube Streaming"
$ws.Range("G1031").Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows(1031).RowHeight = 58

# Row 1032
$ws.Range("A1032").Value = "CH"
$ws.Range("B1032").Value = "MGRN23_CHN_P41"
$ws.Range("C1032").Value = "28-Dec-2025 12:13 PM"
$ws.Range("D1032").Value = "FAIL"
$ws.Range("E1032").Value = "1. RACH setup`n2. SgNB addition Success (ENDC Setup)`n3. MO Call (pass/fail)`n4. MT Call (pass/fail)`n5. Downlink Peak MCS - 5G`n6. Peak Rank - 5G`n7. Ping/Round trip time(ms)`n8. QCI Verification`n9. Video Streaming"
$ws.Range("F1032").Value = "1. Static All`n2. Static ATDT`n3. Static VoLTE MO`n4. Static VoLTE MT`n5. Static DL`n6. Static DL`n7. Static Ping`n8. Static All`n9. Static Yotube Streaming"
$ws.Range("G1032").Value = "1. If DT Tool is  TEMS Pocket, verify the Static ATDT . The NR RACH Attempts should be equal to NR RACH Success; kindly exclude the logs where NR RACH has failed.`nIf DT Tool is  AZQ, Validate Static All and ensure NR RACH Attempts match NR RACH Success. Please exclude the logs with NR RACH failures and redo the test accordingly.`n2. Static ATDT –  Verify that the NR ENDC RRC Reconfiguration count matches the NR ENDC RRC Reconfiguration Complete count. If there is any mismatch, kindly exclude the logfile and create a new one.`n3. VoLTE Long Call MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO calls.`n4. VoLTE Long Call MT – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MT calls.`n5. Peak MCS is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n6. Peak Rank is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n7. Ping is not meeting the acceptance criteria. The average ping value across all logfiles should be less than 50 ms. Kindly exclude the logfile where the average value exceeds 50 ms and redo the test.”`n8. If the DT tool is TEMS, verify the QCI combined value in the DL drive. If the DT tool is AZQ, verify it in Static All. The QCI combined value must meet the acceptance criteria.`n9. While performing the YouTube test for both sectors, please ensure that the video is playing successfully in the script before saving the log file."
$ws.Rows(1032).RowHeight = 409.5

# Row 1033
$ws.Range("A1033").Value = "RJ"
$ws.Range("B1033").Value = "KOT298_JPR_P40"
$ws.Range("C1033").Value = "27-Dec-2025 8:41 PM"
$ws.Range("D1033").Value = "FAIL"
$ws.Range("E1033").Value = "1. CSFB Call (pass/fail)"
$ws.Range("F1033").Value = "1. Static CSFB MO"
$ws.Range("G1033").Value = "1. CSFB MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO  calls."
$ws.Rows(1033).RowHeight = 43.5

# Row 1034
$ws.Range("A1034").Value = "KL"
$ws.Range("B1034").Value = "MNTY06_TSR_P40"
$ws.Range("C1034").Value = "27-Dec-2025 7:55 PM"
$ws.Range("D1034").Value = "FAIL"
$ws.Range("E1034").Value = "1. SCG addition after VoLTE call released`n2. Peak Rank - 5G`n3. SgNB Addition time (ms)"
$ws.Range("F1034").Value = "1. Static VoLTE MO`n2. Static DL`n3. Static ATDT"
$ws.Range("G1034").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n2. Peak Rank is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n3. Exclue ATDT Logfile and Create New Sgnb Addition Time Is Very High. It Should Be <150 Ms. To Achieve This, Perform Static Test In Main Lobe And Keep Test Files Downloading In Background. Also, Ensure 4G Serving Cell Belongs To The Same Site. Exclude The Existing Logfile First"
$ws.Rows(1034).RowHeight = 203

# Row 1035
$ws.Range("A1035").Value = "PB"
$ws.Range("B1035").Value = "LBHGR19_SNR_P40"
$ws.Range("C1035").Value = "28-Dec-2025 9:33 PM"
$ws.Range("D1035").Value = "FAIL"
$ws.Range("E1035").Value = "1. Serving SSB beam steering`n2. SCG addition after VoLTE call released`n3. Ping/Round trip time(ms)"
$ws.Range("F1035").Value = "1. Mobility DL`n2. Static VoLTE MO`n3. Static Ping"
$ws.Range("G1035").Value = "1. Kindly add drive coverage in the failed sector and verify that the Beam Index servings are meeting the acceptance criteria.`n2. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n3. Ping is not meeting the acceptance criteria. The average ping value across all logfiles should be less than 50 ms. Kindly exclude the logfile where the average value exceeds 50 ms and redo the test.”"
$ws.Rows(1035).RowHeight = 159.5

# Row 1036
$ws.Range("A1036").Value = "RJ"
$ws.Range("B1036").Value = "KRSP02_1_JPR_P40"
$ws.Range("C1036").Value = "27-Dec-2025 6:37 PM"
$ws.Range("D1036").Value = "FAIL"
$ws.Range("E1036").Value = "1. CSFB Call (pass/fail)"
$ws.Range("F1036").Value = "1. Static CSFB MO"
$ws.Range("G1036").Value = "1. CSFB MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO  calls."
$ws.Rows(1036).RowHeight = 43.5

# Row 1037
$ws.Range("A1037").Value = "PB"
$ws.Range("B1037").Value = "LTIWN26_CHD_P40"
$ws.Range("C1037").Value = "27-Dec-2025 5:32 PM"
$ws.Range("D1037").Value = "FAIL"
$ws.Range("E1037").Value = "1. SCG addition after VoLTE call released"
$ws.Range("F1037").Value = "1. Static VoLTE MO"
$ws.Range("G1037").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition."
$ws.Rows(1037).RowHeight = 87

# Row 1038
$ws.Range("A1038").Value = "OR"
$ws.Range("B1038").Value = "ORBACH-01_JAJ_P41"
$ws.Range("C1038").Value = "27-Dec-2025 4:35 PM"
$ws.Range("D1038").Value = "FAIL"
$ws.Range("E1038").Value = "1. Video Streaming"
$ws.Range("F1038").Value = "1. Static Yotube Streaming"
$ws.Range("G1038").Value = "1. While performing the YouTube test for both sectors, please ensure that the video is playing successfully in the script before saving the log file."
$ws.Rows(1038).RowHeight = 29

# Row 1039
$ws.Range("A1039").Value = "UW"
$ws.Range("B1039").Value = "AGA219_AGR_P40"
$ws.Range("C1039").Value = "27-Dec-2025 3:02 PM"
$ws.Range("D1039").Value = "FAIL"
$ws.Range("E1039").Value = "1. SCG addition after VoLTE call released`n2. Peak PUSCH UL Throughput"
$ws.Range("F1039").Value = "1. Static VoLTE MO`n2. Static UL"
$ws.Range("G1039").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n2. Peak PUSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PUSCH Throughput in the NR tab."
$ws.Rows(1039).RowHeight = 145

# Row 1040
$ws.Range("A1040").Value = "TN"
$ws.Range("B1040").Value = "CB1221_CBE_P40"
$ws.Range("C1040").Value = "27-Dec-2025 8:22 AM"
$ws.Range("D1040").Value = "FAIL"
$ws.Range("E1040").Value = "1. Video Streaming  (ms)"
$ws.Range("F1040").Value = "1. Static Yotube Streaming"
$ws.Range("G1040").Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows(1040).RowHeight = 58

$ws.Range("E6").Select() | Out-Null
Write-Host "Added rows 1029-1040"